# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
# These two sheets mirror the same exhibition data, so most of the same
# cells are updated on both; "展览" additionally needs its F5 bumped
# since "全部类型" already held the newer value for that row.

$wb = $excel.ActiveWorkbook

$commonUpdates = @{
    "F2"  = 10182
    "F15" = 517
    "F23" = 33
    "F39" = 1563
    "F44" = 3129
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $commonUpdates.Keys) {
        $ws.Range($cellRef).Value = $commonUpdates[$cellRef]
    }
}

# Only the "展览" sheet's F5 changes (57 -> 58); "全部类型" already has 58.
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 58
